# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.850.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.255.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.77%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +4.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.80%  "

$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.588.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.811"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.251.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.724.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.61%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.42%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.20%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.144"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "

$ws.Range("E30").Value = "  +5.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0632"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000240"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.84%  "

$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0240"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0976"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.494.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.82%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.90%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.89%  "
